$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and two reordered rows) per upstream refresh.
$updates = @{
    "D2" = "37.168.79"
    "E2" = "  +0.30%  "
    "D3" = "2.054.08"
    "E3" = "  -0.31%  "
    "E4" = "  +0.21%  "
    "D5" = "248.45"
    "E5" = "  -1.80%  "
    "E6" = "  -0.87%  "
    "E7" = "  -0.03%  "
    "D8" = "56.79"
    "E8" = "  -2.63%  "
    "D9" = "0.385"
    "E9" = "  -0.42%  "
    "D10" = "0.0785"
    "E10" = "  -2.09%  "
    "E11" = "  +0.22%  "
    "D12" = "16.24"
    "E12" = "  -1.62%  "
    "D13" = "0.910"
    "E13" = "  +12.62%  "
    "D14" = "2.351.26"
    "E14" = "  -0.39%  "
    "D15" = "5.77"
    "E15" = "  +2.48%  "
    "D16" = "2.048.68"
    "E16" = "  -0.44%  "
    "D17" = "18.59"
    "E17" = "  +12.77%  "
    "D18" = "37.176.48"
    "D19" = "74.72"
    "E19" = "  -1.36%  "
    "D20" = "0.0₃0901"
    "E20" = "  -2.02%  "
    "D21" = "5.49"
    "E21" = "  +0.19%  "
    "D22" = "237.61"
    "E22" = "  -0.35%  "
    "E23" = "  +0.03%  "
    "D24" = "2.48"
    "E24" = "  +3.27%  "
    "D25" = "9.68"
    "E25" = "  +3.85%  "
    "B26" = "PancakeSwap"
    "C26" = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
    "D26" = "2.18"
    "E26" = "  -5.00%  "
    "B27" = "Monero"
    "C27" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D27" = "170.38"
    "E27" = "  +0.49%  "
    "D28" = "20.24"
    "E28" = "  -0.31%  "
    "E29" = "  -0.76%  "
    "D30" = "5.16"
    "E30" = "  +8.15%  "
    "E31" = "  +1.86%  "
    "D32" = "0.0625"
    "E32" = "  +0.50%  "
    "D33" = "4.57"
    "E33" = "  +1.60%  "
    "D34" = "0.0884"
    "E34" = "  +0.16%  "
    "E35" = "  +0.07%  "
    "E36" = "  -0.52%  "
    "E37" = "  +1.11%  "
    "E38" = "  -1.79%  "
    "D39" = "5.31"
    "E39" = "  +13.75%  "
    "D40" = "3.08"
    "E40" = "  +7.86%  "
    "D41" = "0.100"
    "E41" = "  -11.26%  "
    "B42" = "InjectiveProtocol"
    "C42" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D42" = "17.74"
    "E42" = "  -0.25%  "
    "B43" = "VeChain"
    "C43" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D43" = "0.0224"
    "E43" = "  -0.75%  "
    "E44" = "  -0.46%  "
    "D45" = "96.54"
    "E45" = "  -1.18%  "
    "D46" = "2.46"
    "E46" = "  -1.09%  "
    "D47" = "1.274.25"
    "E48" = "  -1.98%  "
    "D49" = "6.86"
    "E49" = "  -0.86%  "
    "D50" = "2.239.37"
    "E50" = "  -0.47%  "
    "D51" = "44.39"
    "E51" = "  +0.49%  "
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}

